$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows appended at the bottom of the config sheet.
$ws.Range("A24").Value = "sentence_model"

$ws.Range("A25").Value = "sentence_embedding_dims_to_keep"
$ws.Range("B25").Value = 8

# Column widths: widen column A to fit the new, longer keys and give
# column B an explicit width now that it holds a long comma list.
# (inputs chosen so the host's pixel-grid rounding lands on 34 / ~125.875)
$ws.Columns.Item(1).ColumnWidth = 33.3
$ws.Columns.Item(2).ColumnWidth = 125.15

# Scroll/selection state as left by the editor after the change.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
